$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 5.526037
$ws.Range("H2").Value = 16.578111
$ws.Range("I2").Value = 0.1907740523064932
$ws.Range("J2").Value = 0.1907740523064932
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 11.05624166666666
$ws.Range("N2").Value = 33.16872499999999
$ws.Range("O2").Value = 0.2506866656360179
$ws.Range("P2").Value = 0.2506866656360179
$ws.Range("Q2").Value = 61.09720053094166
$ws.Range("R2").Value = 549.8748047784749
$ws.Range("S2").Value = 0.04782451106258606
$ws.Range("T2").Value = 0.04782451106258606
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5.526037
$ws.Range("H3").Value = 16.578111
$ws.Range("I3").Value = 0.1907740523064932
$ws.Range("J3").Value = 0.1907740523064932
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 10.558644
$ws.Range("N3").Value = 31.675932
$ws.Range("O3").Value = 0.2394042512635997
$ws.Range("P3").Value = 0.2394042512635997
$ws.Range("Q3").Value = 58.34745741382799
$ws.Range("R3").Value = 525.127116724452
$ws.Range("S3").Value = 0.04567211915295881
$ws.Range("T3").Value = 0.04567211915295881
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.526037
$ws.Range("H4").Value = 16.578111
$ws.Range("I4").Value = 0.1907740523064932
$ws.Range("J4").Value = 0.1907740523064932
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.822425
$ws.Range("N4").Value = 23.467275
$ws.Range("O4").Value = 0.1773638546948513
$ws.Range("P4").Value = 0.1773638546948513
$ws.Range("Q4").Value = 43.227009979725
$ws.Range("R4").Value = 389.043089817525
$ws.Range("S4").Value = 0.03383642129283684
$ws.Range("T4").Value = 0.03383642129283684
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 5.526037
$ws.Range("H5").Value = 16.578111
$ws.Range("I5").Value = 0.1907740523064932
$ws.Range("J5").Value = 0.1907740523064932
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 14.66651766666667
$ws.Range("N5").Value = 43.999553
$ws.Range("O5").Value = 0.3325452284055311
$ws.Range("P5").Value = 0.3325452284055311
$ws.Range("Q5").Value = 81.04771928715365
$ws.Range("R5").Value = 729.429473584383
$ws.Range("S5").Value = 0.06344100079811153
$ws.Range("T5").Value = 0.06344100079811153
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 8.502533
$ws.Range("H6").Value = 25.507599
$ws.Range("I6").Value = 0.2935309110814287
$ws.Range("J6").Value = 0.2935309110814287
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 11.05624166666666
$ws.Range("N6").Value = 33.16872499999999
$ws.Range("O6").Value = 0.2506866656360179
$ws.Range("P6").Value = 0.2506866656360179
$ws.Range("Q6").Value = 94.00605962680831
$ws.Range("R6").Value = 846.0545366412748
$ws.Range("S6").Value = 0.07358428536010582
$ws.Range("T6").Value = 0.07358428536010582
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 8.502533
$ws.Range("H7").Value = 25.507599
$ws.Range("I7").Value = 0.2935309110814287
$ws.Range("J7").Value = 0.2935309110814287
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 10.558644
$ws.Range("N7").Value = 31.675932
$ws.Range("O7").Value = 0.2394042512635997
$ws.Range("P7").Value = 0.2394042512635997
$ws.Range("Q7").Value = 89.77521904525199
$ws.Range("R7").Value = 807.9769714072679
$ws.Range("S7").Value = 0.07027254799017169
$ws.Range("T7").Value = 0.07027254799017169
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.502533
$ws.Range("H8").Value = 25.507599
$ws.Range("I8").Value = 0.2935309110814287
$ws.Range("J8").Value = 0.2935309110814287
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 7.822425
$ws.Range("N8").Value = 23.467275
$ws.Range("O8").Value = 0.1773638546948513
$ws.Range("P8").Value = 0.1773638546948513
$ws.Range("Q8").Value = 66.510426702525
$ws.Range("R8").Value = 598.593840322725
$ws.Range("S8").Value = 0.05206177386149385
$ws.Range("T8").Value = 0.05206177386149385
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.502533
$ws.Range("H9").Value = 25.507599
$ws.Range("I9").Value = 0.2935309110814287
$ws.Range("J9").Value = 0.2935309110814287
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 14.66651766666667
$ws.Range("N9").Value = 43.999553
$ws.Range("O9").Value = 0.3325452284055311
$ws.Range("P9").Value = 0.3325452284055311
$ws.Range("Q9").Value = 124.7025504559163
$ws.Range("R9").Value = 1122.322954103247
$ws.Range("S9").Value = 0.09761230386965734
$ws.Range("T9").Value = 0.09761230386965734
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 10.61443666666667
$ws.Range("H10").Value = 31.84331
$ws.Range("I10").Value = 0.3664396557335078
$ws.Range("J10").Value = 0.3664396557335078
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 11.05624166666666
$ws.Range("N10").Value = 33.16872499999999
$ws.Range("O10").Value = 0.2506866656360179
$ws.Range("P10").Value = 0.2506866656360179
$ws.Range("Q10").Value = 117.3557769421944
$ws.Range("R10").Value = 1056.20199247975
$ws.Range("S10").Value = 0.09186153545264339
$ws.Range("T10").Value = 0.09186153545264339
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 10.61443666666667
$ws.Range("H11").Value = 31.84331
$ws.Range("I11").Value = 0.3664396557335078
$ws.Range("J11").Value = 0.3664396557335078
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 10.558644
$ws.Range("N11").Value = 31.675932
$ws.Range("O11").Value = 0.2394042512635997
$ws.Range("P11").Value = 0.2394042512635997
$ws.Range("Q11").Value = 112.07405802388
$ws.Range("R11").Value = 1008.66652221492
$ws.Range("S11").Value = 0.08772721141417167
$ws.Range("T11").Value = 0.08772721141417167
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 10.61443666666667
$ws.Range("H12").Value = 31.84331
$ws.Range("I12").Value = 0.3664396557335078
$ws.Range("J12").Value = 0.3664396557335078
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 7.822425
$ws.Range("N12").Value = 23.467275
$ws.Range("O12").Value = 0.1773638546948513
$ws.Range("P12").Value = 0.1773638546948513
$ws.Range("Q12").Value = 83.03063474225
$ws.Range("R12").Value = 747.27571268025
$ws.Range("S12").Value = 0.06499314985394923
$ws.Range("T12").Value = 0.06499314985394923
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 10.61443666666667
$ws.Range("H13").Value = 31.84331
$ws.Range("I13").Value = 0.3664396557335078
$ws.Range("J13").Value = 0.3664396557335078
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 14.66651766666667
$ws.Range("N13").Value = 43.999553
$ws.Range("O13").Value = 0.3325452284055311
$ws.Range("P13").Value = 0.3325452284055311
$ws.Range("Q13").Value = 155.6768228933811
$ws.Range("R13").Value = 1401.09140604043
$ws.Range("S13").Value = 0.1218577590127435
$ws.Range("T13").Value = 0.1218577590127435
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 4.323390666666667
$ws.Range("H14").Value = 12.970172
$ws.Range("I14").Value = 0.1492553808785702
$ws.Range("J14").Value = 0.1492553808785702
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 11.05624166666666
$ws.Range("N14").Value = 33.16872499999999
$ws.Range("O14").Value = 0.2506866656360179
$ws.Range("P14").Value = 0.2506866656360179
$ws.Range("Q14").Value = 47.80045203007778
$ws.Range("R14").Value = 430.2040682707
$ws.Range("S14").Value = 0.03741633376068263
$ws.Range("T14").Value = 0.03741633376068262
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 4.323390666666667
$ws.Range("H15").Value = 12.970172
$ws.Range("I15").Value = 0.1492553808785702
$ws.Range("J15").Value = 0.1492553808785702
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 10.558644
$ws.Range("N15").Value = 31.675932
$ws.Range("O15").Value = 0.2394042512635997
$ws.Range("P15").Value = 0.2394042512635997
$ws.Range("Q15").Value = 45.649142922256
$ws.Range("R15").Value = 410.8422863003041
$ws.Range("S15").Value = 0.0357323727062975
$ws.Range("T15").Value = 0.03573237270629749
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 4.323390666666667
$ws.Range("H16").Value = 12.970172
$ws.Range("I16").Value = 0.1492553808785702
$ws.Range("J16").Value = 0.1492553808785702
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 7.822425
$ws.Range("N16").Value = 23.467275
$ws.Range("O16").Value = 0.1773638546948513
$ws.Range("P16").Value = 0.1773638546948513
$ws.Range("Q16").Value = 33.8193992357
$ws.Range("R16").Value = 304.3745931213001
$ws.Range("S16").Value = 0.02647250968657142
$ws.Range("T16").Value = 0.02647250968657142
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 4.323390666666667
$ws.Range("H17").Value = 12.970172
$ws.Range("I17").Value = 0.1492553808785702
$ws.Range("J17").Value = 0.1492553808785702
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 14.66651766666667
$ws.Range("N17").Value = 43.999553
$ws.Range("O17").Value = 0.3325452284055311
$ws.Range("P17").Value = 0.3325452284055311
$ws.Range("Q17").Value = 63.40908559256845
$ws.Range("R17").Value = 570.6817703331161
$ws.Range("S17").Value = 0.04963416472501868
$ws.Range("T17").Value = 0.04963416472501866
